$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove formulas from column B (B3:B10), replacing with plain literal values.
# Row 5: B5 value changes 3 -> 4, and new I5 = 1
# Row 6: B6 value changes 3.75 -> 4, and new I6 = 0.25
# Row 10: B10 value changes 13.5 -> 20.5; I10 changes 1 -> 11; J10 (=3) removed

$ws.Range("B3").Value = 2.5
$ws.Range("B4").Value = 4
$ws.Range("B5").Value = 4
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 0.5
$ws.Range("B8").Value = 7.5
$ws.Range("B9").Value = 0.5
$ws.Range("B10").Value = 20.5

$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 0.25
$ws.Range("I10").Value = 11

# Remove J10 and the whole row 11 (which held the grand-total formula)
$ws.Range("J10").ClearContents()
$ws.Rows("11").Delete()

# Update selection to match the new target state
$ws.Range("C10").Select()
